$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.024.65'
$ws.Cells.Item(2, 5).Value = '  -3.69%  '
$ws.Cells.Item(3, 4).Value = '1.650.70'
$ws.Cells.Item(3, 5).Value = '  -5.29%  '
$ws.Cells.Item(4, 4).Value = '0.9996'
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).Value = '236.44'
$ws.Cells.Item(5, 5).Value = '  -5.62%  '
$ws.Cells.Item(6, 5).Value = '  -0.02%  '
$ws.Cells.Item(7, 4).Value = '0.4827'
$ws.Cells.Item(7, 5).Value = '  -6.35%  '
$ws.Cells.Item(8, 4).Value = '0.2617'
$ws.Cells.Item(8, 5).Value = '  -5.20%  '
$ws.Cells.Item(9, 4).Value = '0.06018'
$ws.Cells.Item(9, 5).Value = '  -2.89%  '
$ws.Cells.Item(10, 4).Value = '0.07191'
$ws.Cells.Item(10, 5).Value = '  -0.65%  '
$ws.Cells.Item(11, 4).Value = '1.648.61'
$ws.Cells.Item(11, 5).Value = '  -5.43%  '
$ws.Cells.Item(12, 4).Value = '14.78'
$ws.Cells.Item(12, 5).Value = '  -2.77%  '
$ws.Cells.Item(13, 4).Value = '0.6213'
$ws.Cells.Item(13, 5).Value = '  -4.65%  '
$ws.Cells.Item(14, 4).Value = '4.574'
$ws.Cells.Item(14, 5).Value = '  -1.35%  '
$ws.Cells.Item(15, 4).Value = '72.97'
$ws.Cells.Item(15, 5).Value = '  -6.34%  '
$ws.Cells.Item(16, 4).Value = '1.0000'
$ws.Cells.Item(16, 5).Value = '  +0.01%  '
$ws.Cells.Item(17, 4).Value = '0.9999'
$ws.Cells.Item(17, 5).Value = '  +0.00%  '
$ws.Cells.Item(18, 4).Value = '25.009.81'
$ws.Cells.Item(19, 5).Value = '  -3.08%  '
$ws.Cells.Item(20, 4).Value = '0.000006633'
$ws.Cells.Item(20, 5).Value = '  -2.59%  '
$ws.Cells.Item(21, 5).Value = '  +5.76%  '
$ws.Cells.Item(22, 4).Value = '1.858.53'
$ws.Cells.Item(22, 5).Value = '  -5.47%  '
$ws.Cells.Item(23, 4).Value = '8.614'
$ws.Cells.Item(23, 5).Value = '  -0.86%  '
$ws.Cells.Item(24, 4).Value = '5.301'
$ws.Cells.Item(24, 5).Value = '  -1.73%  '
$ws.Cells.Item(25, 4).Value = '132.04'
$ws.Cells.Item(25, 5).Value = '  -2.98%  '
$ws.Cells.Item(26, 5).Value = '  -2.31%  '
$ws.Cells.Item(27, 4).Value = '1.398'
$ws.Cells.Item(27, 5).Value = '  -7.53%  '
$ws.Cells.Item(28, 4).Value = '103.10'
$ws.Cells.Item(28, 5).Value = '  -2.68%  '
$ws.Cells.Item(29, 4).Value = '1.678'
$ws.Cells.Item(29, 5).Value = '  -6.19%  '
$ws.Cells.Item(30, 4).Value = '3.762'
$ws.Cells.Item(30, 5).Value = '  -4.94%  '
$ws.Cells.Item(32, 4).Value = '3.590'
$ws.Cells.Item(32, 5).Value = '  -2.15%  '
$ws.Cells.Item(33, 4).Value = '0.04565'
$ws.Cells.Item(33, 5).Value = '  -2.56%  '
$ws.Cells.Item(34, 2).Value = 'Frax'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(34, 4).Value = '0.9991'
$ws.Cells.Item(34, 5).Value = '  +0.03%  '
$ws.Cells.Item(35, 2).Value = 'HuobiToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(35, 4).Value = '2.592'
$ws.Cells.Item(35, 5).Value = '  -2.43%  '
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(36, 4).Value = '0.9361'
$ws.Cells.Item(36, 5).Value = '  -6.76%  '
$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(37, 4).Value = '0.5795'
$ws.Cells.Item(37, 5).Value = '  -7.47%  '
$ws.Cells.Item(38, 2).Value = 'MXToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(38, 4).Value = '2.597'
$ws.Cells.Item(38, 5).Value = '  -4.94%  '
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39, 4).Value = '0.01563'
$ws.Cells.Item(39, 5).Value = '  -3.36%  '
$ws.Cells.Item(40, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(40, 4).Value = '0.8491'
$ws.Cells.Item(40, 5).Value = '  +11.14%  '
$ws.Cells.Item(41, 2).Value = 'PaxDollar'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(41, 4).Value = '0.9994'
$ws.Cells.Item(41, 5).Value = '  +0.00%  '
$ws.Cells.Item(42, 2).Value = 'RenderToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(42, 4).Value = '1.829'
$ws.Cells.Item(42, 5).Value = '  -4.98%  '
$ws.Cells.Item(43, 2).Value = 'Quant'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(43, 4).Value = '98.36'
$ws.Cells.Item(43, 5).Value = '  -2.14%  '
$ws.Cells.Item(44, 2).Value = 'TheSandbox'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(44, 4).Value = '0.3730'
$ws.Cells.Item(44, 5).Value = '  -3.84%  '
$ws.Cells.Item(45, 2).Value = 'FraxShare'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(45, 4).Value = '4.788'
$ws.Cells.Item(45, 5).Value = '  -4.66%  '
$ws.Cells.Item(46, 2).Value = 'Algorand'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(46, 4).Value = '0.1146'
$ws.Cells.Item(46, 5).Value = '  +0.89%  '
$ws.Cells.Item(47, 2).Value = 'Aptos'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(47, 4).Value = '6.142'
$ws.Cells.Item(47, 5).Value = '  -3.58%  '
$ws.Cells.Item(48, 2).Value = 'Cronos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(48, 4).Value = '0.05195'
$ws.Cells.Item(48, 5).Value = '  -0.65%  '
$ws.Cells.Item(49, 2).Value = 'Elrond'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Cells.Item(49, 4).Value = '29.85'
$ws.Cells.Item(49, 5).Value = '  -3.14%  '
$ws.Cells.Item(50, 2).Value = 'TrueUSD'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Cells.Item(50, 4).Value = '1.001'
$ws.Cells.Item(50, 5).Value = '  -0.15%  '
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).Value = '50.30'
$ws.Cells.Item(51, 5).Value = '  -9.61%  '
